# Update the intensity (time scale) system with per-trigger Standard preset values.
# New master list of time-scale options (shared across all per-trigger rows and the
# CustomTimeScaleProvider entry on the Providers sheet).
$newList = "0.08x | 0.10x | 0.12x | 0.15x | 0.20x | 0.21x | 0.23x | 0.25x | 0.26x | 0.28x | 0.30x | 0.34x | 0.35x | 0.40x | 0.45x | 0.50x"

$wb = $excel.ActiveWorkbook

$menuMock = $wb.Worksheets.Item("Menu Mock")
$providers = $wb.Worksheets.Item("Providers")

# Row -> new Standard value for column D (only rows whose base value actually changes)
$rowUpdates = @{
    43 = "0.28x"   # Basic Kill
    70 = "0.23x"   # Decapitation
    79 = "0.26x"   # Last Enemy
    87 = "0.21x"   # Last Stand
    95 = "0.34x"   # Parry
}

# Every one of these rows' column E (the allowed-values list) gets refreshed,
# including rows 52 (Critical) and 61 (Dismemberment) whose column D value is unchanged.
$rowsToRefreshList = 43, 52, 61, 70, 79, 87, 95

foreach ($row in $rowsToRefreshList) {
    if ($rowUpdates.ContainsKey($row)) {
        $menuMock.Range("D$row").Value = $rowUpdates[$row]
    }
    $menuMock.Range("E$row").Value = $newList
}

# Providers sheet: CustomTimeScaleProvider values list (row 14, column B)
$providers.Range("B14").Value = $newList
